# Auto commit at 2025-08-29  7:46:45.98
# Applies: new custom number format "0.00_ " (style index 2) to the
# charging-volume/revenue columns (C:E), fills in the two newly-completed
# days (rows 56-57) plus blank-but-styled placeholder rows (58-63), and
# updates the sheet's scroll position / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data for the two days that were completed since the last edit ---
$ws.Range("C56").Value = 11819.46
$ws.Range("D56").Value = 10081.41
$ws.Range("E56").Value = 4135.28
$ws.Range("F56").Value = 478

$ws.Range("C57").Value = 4960.5600000000004
$ws.Range("D57").Value = 4170.45
$ws.Range("E57").Value = 1247.9000000000001
$ws.Range("F57").Value = 167

# --- Apply the new "0.00_ " number format to all charging-volume / revenue
#     cells: the already-populated rows (2-57) as well as the still-empty
#     placeholder rows (58-63), which only get the style (no values yet).
#     (Kept as two separate calls - a single multi-area "A,B" range only
#     applies NumberFormat to the first area in this engine.) ---
$ws.Range("C2:E57").NumberFormat = "0.00_ "
$ws.Range("C58:E63").NumberFormat = "0.00_ "

# --- Update the view: scrolled down a bit further, new active cell ---
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("I59").Select()
